# Remove the "język" (language) column from the student-import sheet.
# Before: A..J = klient imie, klient nazwisko, ucznia imie, ucznia nazwisko,
#                plec, JEZYK, grupa, nr tel, email, notatka
# After:  A..I = klient imie, klient nazwisko, ucznia imie, ucznia nazwisko,
#                plec, grupa, nr tel, email, notatka
# (i.e. column F "jezyk" is deleted outright, G..J shift left into F..I)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. stash the exact formatting of the hyperlink cells (I2/I3) before
#        anything moves, so we can restore it after the column shift fixes
#        up the Hyperlinks collection (which otherwise keeps stale refs).
$ws.Range("I2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("I3").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. delete the whole "jezyk" column (F) - this shifts G,H,I,J left by
#        one into F,G,H,I (and drags our temp stash at Z1/Z2 to Y1/Y2 too).
$ws.Columns("F").Delete()

# --- 3. the Hyperlinks collection still anchors at the old I2/I3 location;
#        recreate the two hyperlinks at their new home (H2/H3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", "", "", "olo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ala@gmail.com", "", "", "ala@gmail.com")

# --- 4. re-adding hyperlinks re-styles the cells (built-in Hyperlink
#        style); restore the original manual blue-font formatting.
$ws.Range("Y1").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("Y2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 5. drop the temp stash cells.
$ws.Range("Y1:Y2").Clear()

# --- 6. match the author's final selection (F1, where the removed column
#        used to start).
$ws.Range("F1").Select() | Out-Null
